# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet, with
#    the per-fund holdings table for 600054 (黄山旅游) as of 2022-Q1.
# 2. Update the "总计" (totals) sheet to add a new first data row for
#    2022-Q1 (5 funds, 0.47 亿元 total), re-indexing the existing rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the "2022-Q1" sheet, positioned right before "总计".
# ---------------------------------------------------------------------
$totalsBefore = $wb.Worksheets.Item(4)
$newSheet = $wb.Worksheets.Add($totalsBefore)
$newSheet.Name = "2022-Q1"

# Re-fetch sheet references by index after Add() -- stale references
# captured before the insert silently break Range.Copy() on this host.
$totals = $wb.Worksheets.Item(5)
$q1 = $wb.Worksheets.Item(4)

# Pull the existing bold/bordered/centered header style (cellXfs index 2)
# from the totals sheet and stamp it onto the header row + index column
# of the new sheet, before any values are written.
$totals.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q1.Range("A2:A6").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# index, fund code, fund name, fund size, stock position, position ratio, held value, rank
$fundRows = @(
    @(0, "210003", "金鹰行业优势混合",       "6.62",  "88.02", "3.75", "0.2482", 10),
    @(1, "000968", "广发中证养老产业指数A",   "10.39", "94.08", "1.73", "0.1797", 2),
    @(2, "002982", "广发中证养老产业指数C",   "0.88",  "94.08", "1.73", "0.0152", 2),
    @(3, "516560", "华宝养老ETF",            "0.75",  "97.92", "1.79", "0.0134", 2),
    @(4, "004135", "申万菱信量化成长混合",     "0.49",  "86.91", "1.93", "0.0095", 10)
)

$r = 2
foreach ($row in $fundRows) {
    $q1.Cells.Item($r, 1).Value = $row[0]

    # Fund code / size / position / ratio / value all read as plain text
    # in the source data (leading zeros, e.g. "000968") -- force text so
    # Excel doesn't auto-coerce them to numbers.
    $q1.Cells.Item($r, 2).NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $row[1]

    $q1.Cells.Item($r, 3).Value = $row[2]

    $q1.Cells.Item($r, 4).NumberFormat = "@"
    $q1.Cells.Item($r, 4).Value = $row[3]

    $q1.Cells.Item($r, 5).NumberFormat = "@"
    $q1.Cells.Item($r, 5).Value = $row[4]

    $q1.Cells.Item($r, 6).NumberFormat = "@"
    $q1.Cells.Item($r, 6).Value = $row[5]

    $q1.Cells.Item($r, 7).NumberFormat = "@"
    $q1.Cells.Item($r, 7).Value = $row[6]

    $q1.Cells.Item($r, 8).Value = $row[7]

    # Drop the "@" number-format marker again once the text is committed
    # so the cell is left with no explicit style, matching the rest of
    # the unstyled data cells.
    $q1.Range($q1.Cells.Item($r, 2), $q1.Cells.Item($r, 7)).ClearFormats()

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: rewrite the "总计" sheet data rows, inserting 2022-Q1 on top.
# ---------------------------------------------------------------------
# index, date label, fund count, total held value (亿元)
$totalsRows = @(
    @(0, "2022-Q1", 5, 0.47),
    @(1, "2021-Q3", 2, 0.36),
    @(2, "2021-Q2", 1, 0.03),
    @(3, "2021-Q1", 1, 0.03)
)

$r = 2
foreach ($row in $totalsRows) {
    $totals.Cells.Item($r, 1).Value = $row[0]

    $totals.Cells.Item($r, 2).NumberFormat = "@"
    $totals.Cells.Item($r, 2).Value = $row[1]
    $totals.Cells.Item($r, 2).ClearFormats()

    $totals.Cells.Item($r, 3).Value = $row[2]
    $totals.Cells.Item($r, 4).Value = $row[3]

    $r = $r + 1
}

# Re-stamp the header style (B1:D1, already s=2) onto the whole index
# column -- row 5 is brand new (no inherited style) and row 2's old
# style belonged to what used to be row 4's "2021-Q1" entry, so just
# reapply uniformly across A2:A5 and rewrite the values after.
$totals.Range("B1").Copy()
$totals.Range("A2:A5").PasteSpecial(-4122)
$totals.Cells.Item(2, 1).Value = 0
$totals.Cells.Item(3, 1).Value = 1
$totals.Cells.Item(4, 1).Value = 2
$totals.Cells.Item(5, 1).Value = 3
